$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. First paragraph: append two trailing spaces to the existing text,
#    then append a red "(This is a change - Version for main branch)"
#    comment, typed as three separate chunks so it lands as three runs
#    (matching how the edit actually happened upstream).
# ---------------------------------------------------------------------

function InsertionPointAtEndOf($paraIndex) {
    # Paragraph.Range.End sits just *past* the hidden paragraph mark,
    # so the real "end of visible text" position is End-1. Building a
    # fresh, zero-length Range there (rather than Collapse-ing the
    # paragraph's own Range) is what keeps inserts bound to the live
    # document instead of a detached copy.
    $rng = $d.Paragraphs($paraIndex).Range
    $pos = $rng.End - 1
    return $d.Range($pos, $pos)
}

$ip = InsertionPointAtEndOf 1
$ip.InsertAfter("  ")

$ip = InsertionPointAtEndOf 1
$redStart = $ip.Start
$ip.InsertAfter([char]40 + "This is a change " + [char]8211 + " Ve")
$redEnd = $ip.End
$d.Range($redStart, $redEnd).Font.Color = 255

$ip = InsertionPointAtEndOf 1
$redStart = $ip.Start
$ip.InsertAfter("rsion for main branch")
$redEnd = $ip.End
$d.Range($redStart, $redEnd).Font.Color = 255

$ip = InsertionPointAtEndOf 1
$redStart = $ip.Start
$ip.InsertAfter([char]41)
$redEnd = $ip.End
$d.Range($redStart, $redEnd).Font.Color = 255

Write-Output "Paragraph 1: [$($d.Paragraphs(1).Range.Text)]"

# ---------------------------------------------------------------------
# 2. The blank paragraph right after "It will be treated..." carries
#    explicit Menlo/baseline paragraph + run formatting in its pPr.
#    Strip it down to a bare, totally empty paragraph. Plain property
#    setters (ParagraphFormat/Font) are no-ops on a zero-length range,
#    so instead feed InsertXML a trivial empty-paragraph fragment -
#    that both keeps the paragraph (and its mark) in place and wipes
#    the inherited/explicit formatting in one shot.
# ---------------------------------------------------------------------

$blank = $d.Paragraphs(3).Range
Write-Output "Paragraph 3 before: [$($blank.Text)]"
$blank.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
Write-Output "Paragraph 3 after: [$($d.Paragraphs(3).Range.Text)]"

Write-Output "Total paragraphs: $($d.Paragraphs.Count)"
